$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set "Points completed" column (E) to 0 for rows 2-7
$ws.Range("E2").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("E7").Value = 0

# Update the selection shown in the sheet view
$ws.Range("B13").Select()
